$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string
$ws.Range("A1").Value = "Datos actualizados a 9 de Abril de 2020 a las 12:22"

# Row 5 - España
$ws.Range("B5").Value = 152446
$ws.Range("C5").Value = 4226
$ws.Range("D5").Value = 52165
$ws.Range("E5").Value = 85043
$ws.Range("G5").Value = 446
$ws.Range("H5").Value = 15238

# Row 19 - Austria
$ws.Range("B19").Value = 13028
$ws.Range("C19").Value = 86
$ws.Range("E19").Value = 7493

# Row 32 - Rumania
$ws.Range("B32").Value = 5202
$ws.Range("C32").Value = 441
$ws.Range("D32").Value = 647
$ws.Range("E32").Value = 4326
$ws.Range("F32").Value = 178

# Row 78 - Eslovaquia
$ws.Range("B78").Value = 701
$ws.Range("C78").Value = 19
$ws.Range("E78").Value = 683
